# Auto-generated edit script implementing the diff:
#  - Insert a new worksheet 'Комплайнс хариуцсан нэгжийн үаж' right after 'Project Plan'
#  - Populate it with the new compliance Q&A content (28 new shared strings)
#  - Re-point view/selection state to match the target (tab activation moves to the
#    'Хөдөлмөрийн дотоод журам' sheet, which slides from position 2 to position 3).
$wb = $excel.ActiveWorkbook

# --- Insert new worksheet right after 'Project Plan' (position 2) ---
$wsProjectPlan = $wb.Worksheets.Item(1)
$newWs = $wb.Worksheets.Add($null, $wsProjectPlan)
$newWs.Name = 'Комплайнс хариуцсан нэгжийн үаж'

# Sheet 4 ('МТ-ы үйл ажиллагаа') already carries a bold-title style and a plain-body
# style in the shared style table; copy formats from it so no new duplicate styles
# get created, then tweak the alignment to left (no indent), matching the target.
$wsRef = $wb.Worksheets.Item(4)
$titleFormatCell = $wsRef.Range("A9")
$bodyFormatCell = $wsRef.Range("A11")

$titleStyleSrc = $newWs.Range("B1")
$titleFormatCell.Copy()
$titleStyleSrc.PasteSpecial(-4122)
$titleStyleSrc.HorizontalAlignment = -4131

$bodyStyleSrc = $newWs.Range("B2")
$bodyFormatCell.Copy()
$bodyStyleSrc.PasteSpecial(-4122)
$bodyStyleSrc.IndentLevel = 0
$bodyStyleSrc.HorizontalAlignment = -4131
$bodyStyleSrc.VerticalAlignment = -4108

# --- Write every row, reusing the two formats captured above ---
$c = $titleStyleSrc  # B1 already formatted above
$c.Value = '. Ерөнхий ойлголт'
$newWs.Rows.Item(1).RowHeight = 18

# B2 is the body-format anchor cell; already blank & styled above

$c = $newWs.Range("B3")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'Комплаенс хариуцсан нэгжийн үндсэн зорилго юу вэ?'

$c = $newWs.Range("B4")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B5")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = '“Комплаенс” гэж яг юуг хэлдэг вэ?'

$c = $newWs.Range("B6")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B7")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'Комплаенсын соёл гэж юу вэ, ажилтан бүр яагаад мөрдөх ёстой вэ?'

$c = $newWs.Range("B8")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B9")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'Комплаенсын эрсдэл гэж ямар эрсдэлийг хэлдэг вэ?'

$c = $newWs.Range("B11")
$titleStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = '2. Чиг үүрэг, хамрах хүрээ'
$newWs.Rows.Item(11).RowHeight = 18

$c = $newWs.Range("B12")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B13")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'Комплаенс хариуцсан нэгж ямар чиг үүрэгтэй вэ?'

$c = $newWs.Range("B14")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B15")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'Комплаенсийн бодлого, журамтай нийцүүлэн ажиллахад ажилтнуудын оролцоо ямар байх вэ?'

$c = $newWs.Range("B16")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B17")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'Мөнгө угаах болон терроризмыг санхүүжүүлэхтэй тэмцэх (МУТСТ) үйл ажиллагааг ямар нэгж хариуцдаг вэ?'

$c = $newWs.Range("B18")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B19")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'НҮБ болон олон улсын хориг арга хэмжээг хэрхэн хэрэгжүүлдэг вэ?'

$c = $newWs.Range("B21")
$titleStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = '3. Ажилтны үүрэг, оролцоо'
$newWs.Rows.Item(21).RowHeight = 18

$c = $newWs.Range("B22")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B23")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'Ажилтнууд өдөр тутмын үйл ажиллагаандаа ямар комплаенсын зарчмуудыг мөрдөх ёстой вэ?'

$c = $newWs.Range("B24")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B25")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'Ашиг сонирхлын зөрчил илэрвэл ажилтан яах ёстой вэ?'

$c = $newWs.Range("B26")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B27")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'Харилцагчийн мэдээллийн нууцлалыг хамгаалахад ямар шаардлага тавигддаг вэ?'

$c = $newWs.Range("B28")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B29")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'Зөрчил, дутагдал илрүүлсэн бол хаана мэдээлэх ёстой вэ?'

$c = $newWs.Range("B31")
$titleStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = '4. Эрсдэлийн удирдлага'
$newWs.Rows.Item(31).RowHeight = 18

$c = $newWs.Range("B32")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B33")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'Комплаенсын эрсдэлийг хэрхэн үнэлдэг вэ?'

$c = $newWs.Range("B34")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B35")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'Шинэ бүтээгдэхүүн, үйлчилгээ нэвтрүүлэх үед комплаенсын ямар үнэлгээ, санал шаардлагатай вэ?'

$c = $newWs.Range("B36")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B37")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'Сэжигтэй гүйлгээг илрүүлсэн тохиолдолд яах ёстой вэ?'

$c = $newWs.Range("B38")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B39")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'Өндөр эрсдэлтэй харилцагч болон улс оронтой холбоотой гүйлгээг хэрхэн хянадаг вэ?'

$c = $newWs.Range("B41")
$titleStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = '5. Тайлагнал, хяналт'
$newWs.Rows.Item(41).RowHeight = 18

$c = $newWs.Range("B42")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B43")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'Комплаенсийн нэгж хэнд, хэдэн удаа тайлагнадаг вэ?'

$c = $newWs.Range("B44")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B45")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'МУТСТ-тай холбоотой тайланг хаана, ямар хугацаанд хүргүүлдэг вэ?'

$c = $newWs.Range("B46")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B47")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'Дотоод аудит комплаенсын үйл ажиллагаанд ямар хяналт тавьдаг вэ?'

$c = $newWs.Range("B49")
$titleStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = '6. Хариуцлага'
$newWs.Rows.Item(49).RowHeight = 18

$c = $newWs.Range("B50")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B51")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'Комплаенсын журмыг зөрчвөл ямар хариуцлага хүлээх вэ?'

$c = $newWs.Range("B52")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B53")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'Худал тайлан, буруу мэдээлэл хүргүүлбэл ямар үр дагавартай вэ?'

$c = $newWs.Range("B54")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)

$c = $newWs.Range("B55")
$bodyStyleSrc.Copy()
$c.PasteSpecial(-4122)
$c.Value = 'Хөдөлмөрийн сахилгын шийтгэл хуулийн бусад хариуцлагаас чөлөөлөх үү?'

$excel.CutCopyMode = $false

# --- View/selection state ---
$newWs.Range("B55").Select()

$wsProjectPlanFinal = $wb.Worksheets.Item(1)
$wsProjectPlanFinal.Range("F3").Select()

# 'Хөдөлмөрийн дотоод журам' is now the 3rd tab; make it the active/selected tab
# with its selection moved to A17, matching the target workbook state.
$wsHr = $wb.Worksheets.Item(3)
$wsHr.Activate()
$wsHr.Range("A17").Select()

Write-Host "edit complete"
